$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId=1, sheet1.xml) - rows 2,3,6,9 -> column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 538
$ws1.Range("F3").Value = 6426
$ws1.Range("F6").Value = 137
$ws1.Range("F9").Value = 572

# Sheet "全部类型" (sheetId=4, sheet4.xml) - rows 2,3,7,11 -> column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 538
$ws4.Range("F3").Value = 6426
$ws4.Range("F7").Value = 137
$ws4.Range("F11").Value = 572
